$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers (A1 "Author nm" already present, add the rest)
$ws.Range("B1").Value = "Year"
$ws.Range("C1").Value = "Type"
$ws.Range("D1").Value = "Title"
$ws.Range("E1").Value = "Citation Link"
$ws.Range("F1").Value = "Venue"
$ws.Range("G1").Value = "Publisher"
$ws.Range("H1").Value = "Cited By"
$ws.Range("B1:H1").Style = $ws.Range("A1").Style

# Row 2 - data (fix author name typo, then fill in the rest)
$ws.Range("A2").Value = "Stephan hawking"

$ws.Range("B2").Value = 2019
$ws.Range("C2").Value = "Journal"
$ws.Range("D2").Value = "Security and Privacy Issues in IoT: A Platform for Fog Computing"
$ws.Range("E2").Value = "https://scholar.google.co.in/citations?view_op=view_citation&hl=en&user=xew0uSEAAAAJ&citation_for_view=xew0uSEAAAAJ:2osOgNQ5qMEC"
$ws.Range("F2").Value = "The Journal of Korean Institute of Communications and Information Sciences"
$ws.Range("G2").Value = "N/A"
$ws.Range("H2").Value = 2
$ws.Range("B2:H2").Style = $ws.Range("A1").Style

$ws.Range("G9").Select()
